# feat: add 2022-Q1 data
#
# Inserts a new "2022-Q1" sheet (fund holdings detail) right before the
# existing "总计" (totals) summary sheet, and prepends a matching
# "2022-Q1" row to the "总计" sheet's history table.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 0) Pull off the "总计" sheet's existing history rows, then drop the
#    sheet so it can be recreated (with the new row) at the end of the
#    workbook - this keeps sheet/tab ordering and internal sheet ids in
#    the same relative shape as a plain "insert sheet, then append a
#    refreshed totals sheet" edit.
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")

$history = @()
$r = 2
while ($totalSheet.Cells.Item($r, 2).Value2 -ne $null) {
    $history += , @(
        $totalSheet.Cells.Item($r, 2).Value2,
        $totalSheet.Cells.Item($r, 3).Value2,
        $totalSheet.Cells.Item($r, 4).Value2
    )
    $r = $r + 1
}

$totalSheet.Delete()

# ---------------------------------------------------------------------
# 1) Add the new "2022-Q1" detail sheet at the end of the workbook, and
#    give it the same header/layout as the other quarterly sheets.
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "2022-Q1"

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($c = 0; $c -lt $headers.Length; $c++) {
    $cell = $newSheet.Cells.Item(1, $c + 2)
    $cell.Value = $headers[$c]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
}

$rows = @(
    @("952009", "国泰君安君得鑫两年持有期混合A", "36.44", "80.08", "3.08", "1.1224", 4),
    @("000031", "华夏复兴混合", "27.37", "89.15", "4.04", "1.1057", 10),
    @("952099", "国泰君安君得鑫两年持有期混合C", "28.64", "80.08", "3.08", "0.8821", 4),
    @("160325", "华夏创业板两年定期开放混合", "27.39", "90.77", "2.56", "0.7012", 7),
    @("004139", "中邮军民融合灵活配置混合", "16.83", "86.35", "3.32", "0.5588", 10),
    @("009596", "泰康创新成长混合A", "10.27", "91.52", "3.09", "0.3173", 10),
    @("001479", "中邮风格轮动灵活配置混合", "9.45", "62.17", "3.18", "0.3005", 8),
    @("010106", "华夏核心科技6个月定期开放混合A", "8.53", "79.73", "3.18", "0.2713", 10),
    @("002181", "华安大安全主题灵活配置混合", "5.02", "87.28", "3.19", "0.1601", 6),
    @("160425", "华安创业板两年定期开放混合", "5.11", "96.75", "2.89", "0.1477", 9),
    @("519993", "长信增利动态策略混合", "3.70", "94.82", "3.60", "0.1332", 9),
    @("011927", "博时汇誉回报灵活配置混合型证券投资基金A", "1.30", "68.12", "3.57", "0.0464", 7),
    @("010107", "华夏核心科技6个月定期开放混合C", "1.32", "79.73", "3.18", "0.0420", 10),
    @("002885", "摩根士丹利华鑫万众创新灵活配置混合", "0.86", "94.36", "4.64", "0.0399", 10),
    @("009597", "泰康创新成长混合C", "1.18", "91.52", "3.09", "0.0365", 10),
    @("519971", "长信改革红利灵活配置混合", "0.27", "73.88", "3.86", "0.0104", 6),
    @("011928", "博时汇誉回报灵活配置混合型证券投资基金C", "0.13", "68.12", "3.57", "0.0046", 7)
)

# Text-valued columns (B..G on this sheet) must stay text even though they
# look numeric (codes with leading zeros, fixed-decimal percentages) - so
# force the Text number format before writing the values.
$newSheet.Range("B2:G18").NumberFormat = "@"

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $i + 2
    $data = $rows[$i]

    $idxCell = $newSheet.Cells.Item($r, 1)
    $idxCell.Value = $i
    $idxCell.Font.Bold = $true
    $idxCell.HorizontalAlignment = -4108
    $idxCell.VerticalAlignment = -4160
    $idxCell.Borders.LineStyle = 1

    $newSheet.Cells.Item($r, 2).Value = $data[0]
    $newSheet.Cells.Item($r, 3).Value = $data[1]
    $newSheet.Cells.Item($r, 4).Value = $data[2]
    $newSheet.Cells.Item($r, 5).Value = $data[3]
    $newSheet.Cells.Item($r, 6).Value = $data[4]
    $newSheet.Cells.Item($r, 7).Value = $data[5]
    $newSheet.Cells.Item($r, 8).Value = $data[6]
}

# ---------------------------------------------------------------------
# 2) Recreate "总计" after the new sheet, with a fresh 2022-Q1 row on
#    top of the preserved history.
# ---------------------------------------------------------------------
$newTotal = $wb.Worksheets.Add($null, $newSheet)
$newTotal.Name = "总计"

$totalHeaders = @("日期", "持有数量(只)", "持有市值(亿元)")
for ($c = 0; $c -lt $totalHeaders.Length; $c++) {
    $cell = $newTotal.Cells.Item(1, $c + 2)
    $cell.Value = $totalHeaders[$c]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
}

$totalRows = @(, @("2022-Q1", 17, 5.88)) + $history

for ($i = 0; $i -lt $totalRows.Length; $i++) {
    $r = $i + 2
    $data = $totalRows[$i]

    $idxCell = $newTotal.Cells.Item($r, 1)
    $idxCell.Value = $i
    $idxCell.Font.Bold = $true
    $idxCell.HorizontalAlignment = -4108
    $idxCell.VerticalAlignment = -4160
    $idxCell.Borders.LineStyle = 1

    $newTotal.Cells.Item($r, 2).Value = $data[0]
    $newTotal.Cells.Item($r, 3).Value = $data[1]
    $newTotal.Cells.Item($r, 4).Value = $data[2]
}

# Restore the originally-active sheet/selection so the workbook view isn't
# perturbed by the edits above.
$wb.Worksheets.Item(1).Activate()
$wb.Worksheets.Item(1).Range("A1").Select()
